$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find the BoM row for designator "L1" (the inductor) by scanning column B
# (Designator) rather than hard-coding the row number.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$l1Row = 0
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 2).Value2 -eq "L1") {
        $l1Row = $r
        break
    }
}

if ($l1Row -eq 0) {
    throw "Could not find BoM row for designator L1"
}

# Change L1 inductor value from 1uH to 1.5uH for higher output current
# stability, and update its JLCPCB part number (Comment = column D,
# JLCPCB Part # = column E).
$ws.Cells.Item($l1Row, 4).Value = "1.5uH"
$ws.Cells.Item($l1Row, 5).Value = "C354573"

# Restore the saved selection state.
$ws.Range("E8").Select()
